# Slide 8, shape 6 ("CasellaDiTesto 9") holds the "OUTPUT:" label next to
# the output screenshots. The commit widens the textbox and expands the
# label to call out that the output was tested via Postman.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(8)
$sh = $s.Shapes.Item(6)

# Widen the textbox: cx 3087638 -> 4327436 EMU (PowerPoint COM sizes are
# in points; 12700 EMU per point -> 340.743pt).
$sh.Width = 340.743

# Rebuild the text as three runs: "OUTPUT testato tramite " + "postman" + ":"
$tr = $sh.TextFrame.TextRange
$tr.Text = "OUTPUT testato tramite postman:"

# Touching each Characters() sub-range (even a no-op format re-assert)
# splits the paragraph back into separate <a:r> runs at those boundaries,
# while preserving the existing run formatting (18pt, bg1 solid fill)
# that PowerPoint copies onto each new run.
$r1 = $tr.Characters(1, 23)   # "OUTPUT testato tramite "
$r1.Font.Size = 18

$r2 = $tr.Characters(24, 7)   # "postman"
$r2.Font.Size = 18

$r3 = $tr.Characters(31, 1)   # ":"
$r3.Font.Size = 18
